# implement dgcnn for segmentation
#
# Adds a new "convpoint" sheet, fills in the rest of the npm3d "ours"
# [0:7](v3) row (row 15), adds three new npm3d rows (convpoint, ours
# [0:1](v3), ours [1:2](v3), dgcnn), and swaps the bold/ normal styling
# between rows 14 and 15 to match the new "best row" emphasis.

$wb = $excel.ActiveWorkbook

$npm3d = $wb.Worksheets.Item("npm3d")
$commonClass = $wb.Worksheets.Item("common_class")

# ---------------------------------------------------------------------
# 1. Swap the A14 / A15 styles (bold <-> normal) without fabricating new
#    style entries: round-trip each style through a scratch cell on the
#    SAME sheet (cross-sheet copy/paste loses the font), then delete the
#    scratch row so the used range / dimension isn't polluted.
# ---------------------------------------------------------------------
$a14 = $npm3d.Cells.Item(14, 1)
$a15 = $npm3d.Cells.Item(15, 1)

$scratchRow = 1000
$tempA = $npm3d.Cells.Item($scratchRow, 1)
$tempB = $npm3d.Cells.Item($scratchRow, 2)

$a14.Copy() | Out-Null
$tempA.PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$a15.Copy() | Out-Null
$tempB.PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$tempB.Copy() | Out-Null
$a14.PasteSpecial(-4122) | Out-Null     # A14 <- old A15 style

$tempA.Copy() | Out-Null
$a15.PasteSpecial(-4122) | Out-Null     # A15 <- old A14 style

$npm3d.Rows.Item($scratchRow).Delete() | Out-Null

# ---------------------------------------------------------------------
# 2. Finish row 15 (J15:U15) with the held-out results for that run.
# ---------------------------------------------------------------------
$npm3d.Cells.Item(15, 10).Value = 61
$npm3d.Cells.Item(15, 11).Value = 0.71947099999999997
$npm3d.Cells.Item(15, 12).Value = 0.97270800000000002
$npm3d.Cells.Item(15, 13).Value = 0.97625200000000001
$npm3d.Cells.Item(15, 14).Value = 0.96355299999999999
$npm3d.Cells.Item(15, 15).Value = 0.64958700000000003
$npm3d.Cells.Item(15, 16).Value = 0.62269799999999997
$npm3d.Cells.Item(15, 17).Value = 0.60464399999999996
$npm3d.Cells.Item(15, 18).Value = 0.56358900000000001
$npm3d.Cells.Item(15, 19).Value = 0.29377999999999999
$npm3d.Cells.Item(15, 20).Value = 0.93228800000000001
$npm3d.Cells.Item(15, 21).Value = 0.86884799999999995

# ---------------------------------------------------------------------
# 3. Row 16: convpoint baseline run.
# ---------------------------------------------------------------------
$npm3d.Cells.Item(16, 1).Value = "convpoint"
$a15.Copy() | Out-Null
$npm3d.Cells.Item(16, 1).PasteSpecial(-4122) | Out-Null
$npm3d.Cells.Item(16, 1).Value = "convpoint"

$npm3d.Cells.Item(16, 2).Value = 8192
$npm3d.Cells.Item(16, 3).Value = 24
$npm3d.Cells.Item(16, 4).Value = 1

$npm3d.Cells.Item(14, 5).Copy() | Out-Null
$npm3d.Cells.Item(16, 5).PasteSpecial(-4122) | Out-Null
$npm3d.Cells.Item(16, 5).Value = 0.0001

$npm3d.Cells.Item(16, 6).Value = 0
$npm3d.Cells.Item(16, 7).Value = 0
$npm3d.Cells.Item(16, 8).Value = 10
$npm3d.Cells.Item(16, 9).Value = 10
$npm3d.Cells.Item(16, 10).Value = 70
$npm3d.Cells.Item(16, 11).Value = 0.69466000000000006
$npm3d.Cells.Item(16, 12).Value = 0.97370400000000001
$npm3d.Cells.Item(16, 13).Value = 0.98242300000000005
$npm3d.Cells.Item(16, 14).Value = 0.96644200000000002
$npm3d.Cells.Item(16, 15).Value = 0.67562800000000001
$npm3d.Cells.Item(16, 16).Value = 0.59367899999999996
$npm3d.Cells.Item(16, 17).Value = 0.48457499999999998
$npm3d.Cells.Item(16, 18).Value = 0.54436499999999999
$npm3d.Cells.Item(16, 19).Value = 0.211067
$npm3d.Cells.Item(16, 20).Value = 0.93242199999999997
$npm3d.Cells.Item(16, 21).Value = 0.86133599999999999

# ---------------------------------------------------------------------
# 4. Row 17: another "ours" run, geometry descriptor [0:1](v3).
# ---------------------------------------------------------------------
$npm3d.Cells.Item(17, 1).Value = "ours"
$a15.Copy() | Out-Null
$npm3d.Cells.Item(17, 1).PasteSpecial(-4122) | Out-Null
$npm3d.Cells.Item(17, 1).Value = "ours"

$npm3d.Cells.Item(17, 2).Value = 8192
$npm3d.Cells.Item(17, 3).Value = 16
$npm3d.Cells.Item(17, 4).Value = 1

$npm3d.Cells.Item(14, 5).Copy() | Out-Null
$npm3d.Cells.Item(17, 5).PasteSpecial(-4122) | Out-Null
$npm3d.Cells.Item(17, 5).Value = 0.0001

$npm3d.Cells.Item(17, 6).Value = 0
$npm3d.Cells.Item(17, 7).Value = "[0:1](v3)"
$npm3d.Cells.Item(17, 8).Value = 10
$npm3d.Cells.Item(17, 9).Value = 10
$npm3d.Cells.Item(17, 10).Value = 205
$npm3d.Cells.Item(17, 11).Value = 0.70485200000000003
$npm3d.Cells.Item(17, 12).Value = 0.96953299999999998
$npm3d.Cells.Item(17, 13).Value = 0.97586700000000004
$npm3d.Cells.Item(17, 14).Value = 0.95692999999999995
$npm3d.Cells.Item(17, 15).Value = 0.72497999999999996
$npm3d.Cells.Item(17, 16).Value = 0.60528999999999999
$npm3d.Cells.Item(17, 17).Value = 0.54253799999999996
$npm3d.Cells.Item(17, 18).Value = 0.54025199999999995
$npm3d.Cells.Item(17, 19).Value = 0.26972299999999999
$npm3d.Cells.Item(17, 20).Value = 0.88974500000000001
$npm3d.Cells.Item(17, 21).Value = 0.83834299999999995

# ---------------------------------------------------------------------
# 5. Row 18: "ours" run with geometry descriptor [1:2](v3), still
#    training (no metrics filled in yet).
# ---------------------------------------------------------------------
$npm3d.Cells.Item(18, 1).Value = "ours"
$a15.Copy() | Out-Null
$npm3d.Cells.Item(18, 1).PasteSpecial(-4122) | Out-Null
$npm3d.Cells.Item(18, 1).Value = "ours"

$npm3d.Cells.Item(18, 2).Value = 8192
$npm3d.Cells.Item(18, 3).Value = 16
$npm3d.Cells.Item(18, 4).Value = 1

$npm3d.Cells.Item(14, 5).Copy() | Out-Null
$npm3d.Cells.Item(18, 5).PasteSpecial(-4122) | Out-Null
$npm3d.Cells.Item(18, 5).Value = 0.0001

$npm3d.Cells.Item(18, 6).Value = 0
$npm3d.Cells.Item(18, 7).Value = "[1:2](v3)"
$npm3d.Cells.Item(18, 8).Value = 10
$npm3d.Cells.Item(18, 9).Value = 10

# ---------------------------------------------------------------------
# 6. Row 19: dgcnn baseline, just started (no metrics filled in yet).
# ---------------------------------------------------------------------
$npm3d.Cells.Item(19, 1).Value = "dgcnn"
$a15.Copy() | Out-Null
$npm3d.Cells.Item(19, 1).PasteSpecial(-4122) | Out-Null
$npm3d.Cells.Item(19, 1).Value = "dgcnn"

$npm3d.Cells.Item(19, 2).Value = 8192
$npm3d.Cells.Item(19, 3).Value = 8
$npm3d.Cells.Item(19, 4).Value = 1

$npm3d.Cells.Item(14, 5).Copy() | Out-Null
$npm3d.Cells.Item(19, 5).PasteSpecial(-4122) | Out-Null
$npm3d.Cells.Item(19, 5).Value = 0.0001

$npm3d.Cells.Item(19, 6).Value = 0
$npm3d.Cells.Item(19, 7).Value = 0
$npm3d.Cells.Item(19, 8).Value = 10
$npm3d.Cells.Item(19, 9).Value = 10

# ---------------------------------------------------------------------
# 7. common_class: cursor was left on K1 (no data change).
# ---------------------------------------------------------------------
$commonClass.Select() | Out-Null
$commonClass.Range("K1").Select() | Out-Null

# ---------------------------------------------------------------------
# 8. New "convpoint" sheet, appended after common_class.
# ---------------------------------------------------------------------
$lastIndex = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($lastIndex)
$convpoint = $wb.Worksheets.Add($null, $lastSheet)
$convpoint.Name = "convpoint"

$convpoint.Columns.Item(1).ColumnWidth = 12.75
$convpoint.Columns.Item(2).ColumnWidth = 11.75
$convpoint.Columns.Item(3).ColumnWidth = 11
$convpoint.Columns.Item(4).ColumnWidth = 12.125
$convpoint.Columns.Item(5).ColumnWidth = 13.375
$convpoint.Columns.Item(7).ColumnWidth = 11.375
$convpoint.Columns.Item(8).ColumnWidth = 12.75

$convpoint.Cells.Item(1, 1).Value = "model_name"
$convpoint.Cells.Item(1, 2).Value = "num_point"
$convpoint.Cells.Item(1, 3).Value = "batch_size"
$convpoint.Cells.Item(1, 4).Value = "weight_labels"
$convpoint.Cells.Item(1, 5).Value = "weight_decay"
$convpoint.Cells.Item(1, 6).Value = "use_color"
$convpoint.Cells.Item(1, 7).Value = "use_intensity"
$convpoint.Cells.Item(1, 8).Value = "use_geometry"
$convpoint.Cells.Item(1, 9).Value = "box_size_x"
$convpoint.Cells.Item(1, 10).Value = "box_size_y"
$convpoint.Cells.Item(1, 11).Value = "epoch"
$convpoint.Cells.Item(1, 12).Value = "mIoU"
$convpoint.Cells.Item(1, 13).Value = "OA"

$convpoint.Cells.Item(2, 1).Value = "convpoint"
$convpoint.Cells.Item(2, 2).Value = 8192
$convpoint.Cells.Item(2, 3).Value = 16
$convpoint.Cells.Item(2, 4).Value = 0
$convpoint.Cells.Item(2, 5).Value = 0
$convpoint.Cells.Item(2, 6).Value = 0
$convpoint.Cells.Item(2, 7).Value = 1
$convpoint.Cells.Item(2, 8).Value = 0
$convpoint.Cells.Item(2, 9).Value = 8
$convpoint.Cells.Item(2, 10).Value = 8

$convpoint.Range("H14").Select() | Out-Null

# ---------------------------------------------------------------------
# 9. Leave npm3d as the active sheet/selection (M18), matching the
#    author's last cursor position.
# ---------------------------------------------------------------------
$npm3d.Select() | Out-Null
$npm3d.Range("M18").Select() | Out-Null
